# Auto-generated Excel COM-interop script
# Rewrites rows 2-41 of the active worksheet to the post-edit values
# described by the target diff (new 20-version stimulus set, block_total 3 -> 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order matches the sheet header: A..V
$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

$rowData = @{
    2 = @(93, "categorization", 6, 1, 1, 325, "living_rooms", "living_rooms", "target", $null, "j", "stimuli/img_wgddx.png", 45.6304347826087, 34.30434782608695, 39.96739130434783, 46, 3, 3, 3, 3, 3, 4)
    3 = @(93, "categorization", 6, 1, 2, 326, "living_rooms", "living_rooms", "target", $null, "j", "stimuli/img_a9he3.png", 83.06521739130434, 63.95652173913044, 73.51086956521739, 46, 8, 8, 8, 8, 8, 8)
    4 = @(93, "categorization", 6, 1, 3, 327, "living_rooms", "living_rooms", "target", $null, "j", "stimuli/img_pbsj1.png", 73.88636363636364, 51.52272727272727, 62.70454545454545, 44, 6, 6, 6, 6, 6, 6)
    5 = @(93, "categorization", 6, 1, 4, 328, "living_rooms", "living_rooms", "target", $null, "j", "stimuli/img_i6wsx.png", 79.07142857142857, 58, 68.53571428571428, 42, 7, 7, 7, 7, 7, 7)
    6 = @(93, "categorization", 6, 1, 5, 329, "living_rooms", "bedrooms", "distractor", $null, "f", "stimuli/img_ymgcb.png", 83.73684210526316, 61.13157894736842, 72.4342105263158, 38, 8, 8, 8, 8, 8, 8)
    7 = @(93, "categorization", 6, 1, 6, 330, "living_rooms", "bedrooms", "distractor", $null, "f", "stimuli/img_7caxh.png", 83.38095238095238, 64.26190476190476, 73.82142857142857, 42, 8, 8, 8, 8, 8, 8)
    8 = @(93, "categorization", 6, 1, 7, 331, "living_rooms", "living_rooms", "target", $null, "j", "stimuli/img_dg5h7.png", 88.72093023255815, 76.06976744186046, 82.3953488372093, 43, 10, 10, 10, 10, 10, 10)
    9 = @(93, "categorization", 6, 1, 8, 332, "living_rooms", "kitchens", "distractor", $null, "f", "stimuli/img_xtb79.png", 55.59375, 31.40625, 43.5, 32, 2, 2, 2, 2, 2, 2)
    10 = @(93, "categorization", 6, 1, 9, 333, "living_rooms", "living_rooms", "target", $null, "j", "stimuli/img_w8yhd.png", 55.74418604651163, 38.90697674418605, 47.32558139534883, 43, 4, 4, 4, 4, 4, 4)
    11 = @(93, "categorization", 6, 1, 10, 334, "living_rooms", "living_rooms", "target", $null, "j", "stimuli/img_rych7.png", 30.4468085106383, 23.4468085106383, 26.9468085106383, 47, 2, 2, 2, 2, 2, 2)
    12 = @(93, "categorization", 6, 1, 11, 335, "living_rooms", "living_rooms", "target", $null, "j", "stimuli/img_nb8p4.png", 16.36170212765957, 12.70212765957447, 14.53191489361702, 47, 1, 1, 1, 1, 1, 1)
    13 = @(93, "categorization", 6, 1, 12, 336, "living_rooms", "living_rooms", "target", $null, "j", "stimuli/img_jkm86.png", 58.32558139534883, 38.65116279069768, 48.48837209302326, 43, 4, 4, 4, 4, 4, 4)
    14 = @(93, "categorization", 6, 1, 13, 337, "living_rooms", "living_rooms", "target", $null, "j", "stimuli/img_3sw8t.png", 67.4888888888889, 48.51111111111111, 58, 45, 5, 5, 5, 5, 5, 5)
    15 = @(93, "categorization", 6, 1, 14, 338, "living_rooms", "bedrooms", "distractor", $null, "f", "stimuli/img_swcci.png", 49.82926829268293, 28.46341463414634, 39.14634146341464, 41, 2, 2, 2, 2, 2, 2)
    16 = @(93, "categorization", 6, 1, 15, 339, "living_rooms", "living_rooms", "target", $null, "j", "stimuli/img_xzyzy.png", 85.37209302325581, 68.90697674418605, 77.13953488372093, 43, 9, 9, 9, 9, 9, 9)
    17 = @(93, "categorization", 6, 1, 16, 340, "living_rooms", "kitchens", "distractor", $null, "f", "stimuli/img_7wquy.png", 50.59375, 30.59375, 40.59375, 32, 2, 2, 2, 2, 2, 2)
    18 = @(93, "categorization", 6, 1, 17, 341, "living_rooms", "bedrooms", "distractor", $null, "f", "stimuli/img_ys3qz.png", 46.79545454545455, 31.20454545454545, 39, 44, 2, 2, 2, 2, 2, 2)
    19 = @(93, "categorization", 6, 1, 18, 342, "living_rooms", "kitchens", "distractor", $null, "f", "stimuli/img_3tnh4.png", 80.43243243243244, 58.72972972972973, 69.58108108108108, 37, 8, 8, 8, 8, 8, 8)
    20 = @(93, "categorization", 6, 1, 19, 343, "living_rooms", "kitchens", "distractor", $null, "f", "stimuli/img_gztbt.png", 55.06451612903226, 26.09677419354839, 40.58064516129032, 31, 2, 2, 2, 2, 2, 2)
    21 = @(93, "categorization", 6, 1, 20, 344, "living_rooms", "bedrooms", "distractor", $null, "f", "stimuli/img_b2jli.png", 83.15625, 63.8125, 73.484375, 32, 8, 8, 8, 8, 8, 8)
    22 = @(93, "categorization", 6, 1, 21, 345, "living_rooms", "kitchens", "distractor", $null, "f", "stimuli/img_ub9nn.png", 78.77142857142857, 60.37142857142857, 69.57142857142857, 35, 8, 8, 8, 8, 8, 8)
    23 = @(93, "categorization", 6, 1, 22, 346, "living_rooms", "kitchens", "distractor", $null, "f", "stimuli/img_463mq.png", 51.35294117647059, 30.20588235294118, 40.77941176470588, 34, 2, 2, 2, 2, 2, 2)
    24 = @(93, "categorization", 6, 1, 23, 347, "living_rooms", "kitchens", "distractor", $null, "f", "stimuli/img_uegbb.png", 78.80952380952381, 61.52380952380953, 70.16666666666667, 42, 8, 8, 8, 8, 8, 8)
    25 = @(93, "categorization", 6, 1, 24, 348, "living_rooms", "bedrooms", "distractor", $null, "f", "stimuli/img_91csq.png", 50.44736842105263, 28.34210526315789, 39.39473684210526, 38, 2, 2, 2, 2, 2, 2)
    26 = @(93, "categorization", 6, 1, 25, 349, "living_rooms", "living_rooms", "target", $null, "j", "stimuli/img_o30wb.png", 81.06666666666666, 65.37777777777778, 73.22222222222223, 45, 8, 8, 8, 8, 8, 8)
    27 = @(93, "categorization", 6, 1, 26, 350, "living_rooms", "living_rooms", "target", $null, "j", "stimuli/img_8dmpq.png", 30.65909090909091, 24.11363636363636, 27.38636363636364, 44, 2, 2, 2, 2, 2, 2)
    28 = @(93, "categorization", 6, 1, 27, 351, "living_rooms", "bedrooms", "distractor", $null, "f", "stimuli/img_th7xh.png", 82.35897435897436, 65.53846153846153, 73.94871794871796, 39, 8, 8, 8, 8, 8, 8)
    29 = @(93, "categorization", 6, 1, 28, 352, "living_rooms", "kitchens", "distractor", $null, "f", "stimuli/img_c30d1.png", 78.875, 60.34375, 69.609375, 32, 8, 8, 8, 8, 8, 8)
    30 = @(93, "categorization", 6, 1, 29, 353, "living_rooms", "living_rooms", "target", $null, "j", "stimuli/img_g13d5.png", 73, 51.51111111111111, 62.25555555555556, 45, 6, 6, 6, 6, 6, 6)
    31 = @(93, "categorization", 6, 1, 30, 354, "living_rooms", "living_rooms", "target", $null, "j", "stimuli/img_c0vzo.png", 21.51162790697675, 8.232558139534884, 14.87209302325581, 43, 1, 1, 1, 1, 1, 1)
    32 = @(93, "categorization", 6, 1, 31, 355, "living_rooms", "living_rooms", "target", $null, "j", "stimuli/img_zxvl3.png", 68.78260869565217, 47.56521739130435, 58.17391304347827, 46, 5, 5, 5, 5, 5, 5)
    33 = @(93, "categorization", 6, 1, 32, 356, "living_rooms", "kitchens", "distractor", $null, "f", "stimuli/img_uhmld.png", 54.19354838709678, 32, 43.09677419354838, 31, 2, 2, 2, 2, 2, 2)
    34 = @(93, "categorization", 6, 1, 33, 357, "living_rooms", "living_rooms", "target", $null, "j", "stimuli/img_jpldg.png", 79.54545454545455, 57.75, 68.64772727272728, 44, 7, 7, 7, 7, 7, 7)
    35 = @(93, "categorization", 6, 1, 34, 358, "living_rooms", "bedrooms", "distractor", $null, "f", "stimuli/img_gv750.png", 83.51428571428572, 61.88571428571429, 72.7, 35, 8, 8, 8, 8, 8, 8)
    36 = @(93, "categorization", 6, 1, 35, 359, "living_rooms", "living_rooms", "target", $null, "j", "stimuli/img_165pk.png", 85.73333333333333, 69.22222222222223, 77.47777777777779, 45, 9, 9, 9, 9, 9, 9)
    37 = @(93, "categorization", 6, 1, 36, 360, "living_rooms", "living_rooms", "target", $null, "j", "stimuli/img_5jy9c.png", 87.37209302325581, 79.18604651162791, 83.27906976744185, 43, 10, 10, 10, 10, 9, 10)
    38 = @(93, "categorization", 6, 1, 37, 361, "living_rooms", "kitchens", "distractor", $null, "f", "stimuli/img_lgyo6.png", 54.68421052631579, 27.18421052631579, 40.93421052631579, 38, 2, 2, 2, 2, 2, 2)
    39 = @(93, "categorization", 6, 1, 38, 362, "living_rooms", "bedrooms", "distractor", $null, "f", "stimuli/img_awuev.png", 44.21052631578947, 24.26315789473684, 34.23684210526316, 38, 2, 2, 2, 2, 2, 2)
    40 = @(93, "categorization", 6, 1, 39, 363, "living_rooms", "living_rooms", "target", $null, "j", "stimuli/img_9bkl9.png", 46.62162162162162, 34.27027027027027, 40.44594594594595, 37, 3, 3, 3, 3, 3, 3)
    41 = @(93, "categorization", 6, 1, 40, 364, "living_rooms", "bedrooms", "distractor", $null, "f", "stimuli/img_md6k3.png", 81.57142857142857, 63.66666666666666, 72.61904761904762, 42, 8, 8, 8, 8, 8, 8)
}

foreach ($r in $rowData.Keys) {
    $values = $rowData[$r]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $col = $columns[$i]
        $val = $values[$i]
        if ($null -ne $val) {
            $ws.Range("$col$r").Value = $val
        }
    }
}

Write-Host "Applied living_rooms categorization refresh to rows 2-41."
